$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells keep their original text (not numeric) representation,
# matching the source data which stores these as inline strings (e.g. "60.844.69", "7.05").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.844.69"
$ws.Range("E2").Value = "  +3.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.709.12"
$ws.Range("E3").Value = "  +2.71%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.57"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.31"
$ws.Range("E6").Value = "  +2.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.579"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.737.54"
$ws.Range("E9").Value = "  +3.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.05"
$ws.Range("E10").Value = "  +10.85%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.343"
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.187.34"
$ws.Range("E14").Value = "  +2.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.812.90"
$ws.Range("E15").Value = "  +3.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.61"
$ws.Range("E16").Value = "  +3.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.733.54"
$ws.Range("E17").Value = "  +3.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000139"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "346.35"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.60"
$ws.Range("E21").Value = "  +3.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.41"
$ws.Range("E22").Value = "  +3.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.77"
$ws.Range("E24").Value = "  +2.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.172"
$ws.Range("E25").Value = "  +5.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.421"
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0830"
$ws.Range("E28").Value = "  +2.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.33"
$ws.Range("E29").Value = "  +3.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.78"
$ws.Range("E30").Value = "  +7.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.60"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.16"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "151.33"
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.30"
$ws.Range("E35").Value = "  +6.77%  "
$ws.Range("E36").Value = "  +6.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.924"
$ws.Range("E37").Value = "  -5.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.909"
$ws.Range("E38").Value = "  +7.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.54"
$ws.Range("E39").Value = "  +7.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.51"
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.68"
$ws.Range("E41").Value = "  +0.85%  "
$ws.Range("E42").Value = "  +4.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "282.87"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.19"
$ws.Range("E44").Value = "  +2.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0988"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.98"
$ws.Range("E47").Value = "  +5.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0547"
$ws.Range("E48").Value = "  +4.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.094.84"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.61"
$ws.Range("E50").Value = "  +4.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.54"
$ws.Range("E51").Value = "  +2.22%  "
